$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.854221333333333
$ws.Range("H2").Value = 5.562664
$ws.Range("I2").Value = 0.03416002559055492
$ws.Range("J2").Value = 0.03416002559055492
$ws.Range("M2").Value = 123.2806423333333
$ws.Range("N2").Value = 369.841927
$ws.Range("O2").Value = 0.6241574062367528
$ws.Range("P2").Value = 0.6241574062367526
$ws.Range("Q2").Value = 228.5895970015031
$ws.Range("R2").Value = 2057.306373013528
$ws.Range("S2").Value = 0.02132123296958186
$ws.Range("T2").Value = 0.02132123296958185
$ws.Range("G3").Value = 1.854221333333333
$ws.Range("H3").Value = 5.562664
$ws.Range("I3").Value = 0.03416002559055492
$ws.Range("J3").Value = 0.03416002559055492
$ws.Range("O3").Value = 0.2392728888301323
$ws.Range("P3").Value = 0.2392728888301322
$ws.Range("Q3").Value = 87.63060837624442
$ws.Range("R3").Value = 788.6754753861999
$ws.Range("S3").Value = 0.008173568005563321
$ws.Range("T3").Value = 0.00817356800556332
$ws.Range("G4").Value = 1.854221333333333
$ws.Range("H4").Value = 5.562664
$ws.Range("I4").Value = 0.03416002559055492
$ws.Range("J4").Value = 0.03416002559055492
$ws.Range("O4").Value = 0.136569704933115
$ws.Range("P4").Value = 0.136569704933115
$ws.Range("Q4").Value = 50.01689237575645
$ws.Range("R4").Value = 450.152031381808
$ws.Range("S4").Value = 0.004665224615409745
$ws.Range("T4").Value = 0.004665224615409744
$ws.Range("I5").Value = 0.8311547934421808
$ws.Range("J5").Value = 0.8311547934421808
$ws.Range("M5").Value = 123.2806423333333
$ws.Range("N5").Value = 369.841927
$ws.Range("O5").Value = 0.6241574062367528
$ws.Range("P5").Value = 0.6241574062367526
$ws.Range("Q5").Value = 5561.861737344479
$ws.Range("R5").Value = 50056.75563610032
$ws.Range("S5").Value = 0.5187714200561155
$ws.Range("T5").Value = 0.5187714200561154
$ws.Range("I6").Value = 0.8311547934421808
$ws.Range("J6").Value = 0.8311547934421808
$ws.Range("O6").Value = 0.2392728888301323
$ws.Range("P6").Value = 0.2392728888301322
$ws.Range("S6").Value = 0.1988728084919225
$ws.Range("T6").Value = 0.1988728084919224
$ws.Range("I7").Value = 0.8311547934421808
$ws.Range("J7").Value = 0.8311547934421808
$ws.Range("O7").Value = 0.136569704933115
$ws.Range("P7").Value = 0.136569704933115
$ws.Range("S7").Value = 0.1135105648941428
$ws.Range("T7").Value = 0.1135105648941428
$ws.Range("I8").Value = 0.1346851809672642
$ws.Range("J8").Value = 0.1346851809672642
$ws.Range("M8").Value = 123.2806423333333
$ws.Range("N8").Value = 369.841927
$ws.Range("O8").Value = 0.6241574062367528
$ws.Range("P8").Value = 0.6241574062367526
$ws.Range("Q8").Value = 901.276585925453
$ws.Range("R8").Value = 8111.489273329077
$ws.Range("S8").Value = 0.0840647532110553
$ws.Range("T8").Value = 0.08406475321105529
$ws.Range("I9").Value = 0.1346851809672642
$ws.Range("J9").Value = 0.1346851809672642
$ws.Range("O9").Value = 0.2392728888301323
$ws.Range("P9").Value = 0.2392728888301322
$ws.Range("S9").Value = 0.03222651233264646
$ws.Range("T9").Value = 0.03222651233264646
$ws.Range("I10").Value = 0.1346851809672642
$ws.Range("J10").Value = 0.1346851809672642
$ws.Range("O10").Value = 0.136569704933115
$ws.Range("P10").Value = 0.136569704933115
$ws.Range("S10").Value = 0.01839391542356248
$ws.Range("T10").Value = 0.01839391542356247
